$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells that already exist in the sheet (empty, but with the right style) ---
# Just set their value; the existing style on the cell is preserved.
$ws.Range("C6").Value = "Bëlga - Rendőrmunka"
$ws.Range("C11").Value = "Ivan & The Parazol - Take My Hand"
$ws.Range("B12").Value = "Sipi"
$ws.Range("C12").Value = "Groovehouse - Vándor"
$ws.Range("G12").Value = "No information available"

# --- Cells that do not exist yet: need to create with matching style ---
# Copy formatting from a sibling cell in the same row/column family first.
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "No information available"

$ws.Range("F6").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = "No information available"

$ws.Range("F9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = "No information available"

$ws.Range("F10").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("G10").Value = "No information available"

$ws.Range("F11").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Guli"

$ws.Range("F11").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("G11").Value = "No information available"

# --- Plain value updates on cells that already exist with content/style ---
$ws.Range("G5").Value = "20g"
$ws.Range("C8").Value = "Bee Gees - Stayin Alive "

# Walk-On Music column got wider now that every row has an entry (bestFit recalculated by Excel)
$ws.Columns("C:C").ColumnWidth = 28.33

$ws.Range("B14").Select()
